$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two existing values that changed
$ws.Range("B2").Value = "HDFC659df"
$ws.Range("B3").Value = "AXIS01"

# Add the new "Description" column (C) with its four values
$ws.Range("C1:C4").NumberFormat = "@"
$ws.Range("C1").Value = "Description"
$ws.Range("C2").Value = "Banglore"
$ws.Range("C3").Value = "Delhi"
$ws.Range("C4").Value = "Jaipur"

# Give column C an explicit custom width of 15 (Excel's ColumnWidth property
# reports ~5/6 of a character less than the stored OOXML column width)
$ws.Columns.Item(3).ColumnWidth = 14.1667

# Move the active selection to B3, matching the saved view state
$ws.Range("B3").Select()

# Resize the workbook window to match the saved view state
$win = $wb.Windows.Item(1)
$win.Width = 12830
$win.Height = 3000
